$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CV (address) value from 6019 to 6018; dependent formula in B4 recalculates automatically
$ws.Range("B1").Value = 6018

# Update selection to B2 on the active sheet
$ws.Range("B2").Select()
